$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, "YCHERN", "ASFLI", "DEREGISTERPROJECT", "PENDING", 1, "", ""),
    @(2, "YCHERN", "ASFLI", "DEREGISTERPROJECT", "PENDING", 1, "", ""),
    @(3, "YCHERN", "", "CHANGETITLE", "PENDING", 0, "test", ""),
    @(4, "YCHERN", "ASMADHUKUMAR", "CHANGETITLE", "PENDING", 1, "yays", "")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 1
        $value = $rowData[$j]
        $cell = $ws.Cells.Item($row, $col)
        if ($value -eq "") {
            # Materialize the blank cell (with default style) without
            # actually storing a value, matching a touched-but-empty cell.
            $cell.Font.Bold = $false
        } else {
            $cell.Value = $value
        }
    }
}

$ws.Range("H5").Select()
